$d = $word.ActiveDocument

# --- Update the date heading ---
$d.Content.Find.Execute("2024-08-09 Friday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-08-10 Saturday", 2)

# --- Update the multiplication problems in the table ---
# Each entry maps (row, column) -> new text. Using direct cell addressing
# (rather than Find/Replace) avoids any ambiguity from new values that
# coincide with other original values elsewhere in the table.
$t = $d.Tables(1)

$updates = @(
    @{ Row = 1;  Col = 1; Text = "652×4=" },
    @{ Row = 1;  Col = 2; Text = "860×8=" },
    @{ Row = 1;  Col = 3; Text = "694×4=" },
    @{ Row = 1;  Col = 4; Text = "220×8=" },
    @{ Row = 1;  Col = 5; Text = "864×5=" },

    @{ Row = 5;  Col = 1; Text = "511×8=" },
    @{ Row = 5;  Col = 2; Text = "290×4=" },
    @{ Row = 5;  Col = 3; Text = "147×7=" },
    @{ Row = 5;  Col = 4; Text = "300×8=" },
    @{ Row = 5;  Col = 5; Text = "523×4=" },

    @{ Row = 10; Col = 1; Text = "446×7=" },
    @{ Row = 10; Col = 2; Text = "963×8=" },
    @{ Row = 10; Col = 3; Text = "775×7=" },
    @{ Row = 10; Col = 4; Text = "299×7=" },
    @{ Row = 10; Col = 5; Text = "107×8=" },

    @{ Row = 15; Col = 1; Text = "285×7=" },
    @{ Row = 15; Col = 2; Text = "434×9=" },
    @{ Row = 15; Col = 3; Text = "262×2=" },
    @{ Row = 15; Col = 4; Text = "451×2=" },
    @{ Row = 15; Col = 5; Text = "283×4=" },

    @{ Row = 20; Col = 1; Text = "425×6=" },
    @{ Row = 20; Col = 2; Text = "578×3=" },
    @{ Row = 20; Col = 3; Text = "404×6=" },
    @{ Row = 20; Col = 4; Text = "576×4=" },
    @{ Row = 20; Col = 5; Text = "729×5=" }
)

foreach ($u in $updates) {
    $t.Cell($u.Row, $u.Col).Range.Text = $u.Text
}
